$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 241.625
$ws.Range("I33").Value = 228.91667
$ws.Range("J33").Value = 279.75
$ws.Range("K33").Value = 228.91667
$ws.Range("L33").Value = 279.75
$ws.Range("M33").Value = 0.08332999999998947
$ws.Range("N33").Value = -737.75
$ws.Range("H92").Value = 1514.2142
$ws.Range("I92").Value = 1450.875
$ws.Range("J92").Value = 1598.6666
$ws.Range("K92").Value = 1450.875
$ws.Range("L92").Value = 1598.6666
$ws.Range("M92").Value = -202.875
$ws.Range("N92").Value = -4094.6666
$ws.Range("H93").Value = 34566.086
$ws.Range("J93").Value = 34566.086
$ws.Range("L93").Value = 34566.086
$ws.Range("N93").Value = -39558.086
$ws.Range("H112").Value = 27028440
$ws.Range("I112").Value = 333333800
$ws.Range("K112").Value = 1000001400
$ws.Range("M112").Value = -1000000292
$ws.Range("H129").Value = 819.8200000000001
$ws.Range("I129").Value = 331.7143
$ws.Range("K129").Value = 995.1428999999999
$ws.Range("M129").Value = 4004.8571
$ws.Range("H137").Value = 1324777.5
$ws.Range("I137").Value = 1985476.9
$ws.Range("J137").Value = 3378.8333
$ws.Range("K137").Value = 5956430.699999999
$ws.Range("L137").Value = 10136.4999
$ws.Range("M137").Value = -5953880.699999999
$ws.Range("N137").Value = -15236.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1251.8889
$ws.Range("I2").Value = 1251.8889
$ws.Range("K2").Value = 1251.8889
$ws.Range("M2").Value = -1138.8889
$ws.Range("H32").Value = 4086.0676
$ws.Range("I32").Value = 3410.6064
$ws.Range("J32").Value = 7255.5386
$ws.Range("K32").Value = 3410.6064
$ws.Range("L32").Value = 7255.5386
$ws.Range("M32").Value = -3123.6064
$ws.Range("N32").Value = -7829.5386
$ws.Range("H74").Value = 4539.241
$ws.Range("I74").Value = 5065.579
$ws.Range("K74").Value = 5065.579
$ws.Range("M74").Value = -4191.579
$ws.Range("H77").Value = 4539.241
$ws.Range("I77").Value = 5065.579
$ws.Range("K77").Value = 25327.895
$ws.Range("M77").Value = -20959.895
$ws.Range("H103").Value = 34710
$ws.Range("J103").Value = 34710
$ws.Range("L103").Value = 34710
$ws.Range("N103").Value = -37054
$ws.Range("H116").Value = 1251.8889
$ws.Range("I116").Value = 1251.8889
$ws.Range("K116").Value = 1251.8889
$ws.Range("M116").Value = 1042.1111
$ws.Range("H132").Value = 1772.1526
$ws.Range("I132").Value = 1197.2162
$ws.Range("J132").Value = 2739.0908
$ws.Range("K132").Value = 3591.6486
$ws.Range("L132").Value = 8217.2724
$ws.Range("M132").Value = -1061.6486
$ws.Range("N132").Value = -13277.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1251.8889
$ws.Range("I3").Value = 1251.8889
$ws.Range("K3").Value = 1251.8889
$ws.Range("M3").Value = -1137.8889
$ws.Range("H134").Value = 2767.92
$ws.Range("I134").Value = 1044.6842
$ws.Range("J134").Value = 8224.833000000001
$ws.Range("K134").Value = 3134.0526
$ws.Range("L134").Value = 24674.499
$ws.Range("M134").Value = -599.0526
$ws.Range("N134").Value = -29744.499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2760.0938
$ws.Range("I31").Value = 1018.86365
$ws.Range("K31").Value = 1018.86365
$ws.Range("M31").Value = -723.86365
$ws.Range("H34").Value = 2760.0938
$ws.Range("I34").Value = 1018.86365
$ws.Range("K34").Value = 1018.86365
$ws.Range("M34").Value = -816.86365
$ws.Range("H58").Value = 2566.6086
$ws.Range("I58").Value = 1697.4067
$ws.Range("K58").Value = 1697.4067
$ws.Range("M58").Value = -1494.4067
$ws.Range("H134").Value = 1484.425
$ws.Range("I134").Value = 873.36
$ws.Range("J134").Value = 2502.8667
$ws.Range("K134").Value = 2620.08
$ws.Range("L134").Value = 7508.6001
$ws.Range("M134").Value = -85.07999999999993
$ws.Range("N134").Value = -12578.6001
$ws.Range("H136").Value = 2566.6086
$ws.Range("I136").Value = 1697.4067
$ws.Range("K136").Value = 5092.2201
$ws.Range("M136").Value = -2542.2201

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1404.9474
$ws.Range("I5").Value = 441
$ws.Range("J5").Value = 2272.5
$ws.Range("K5").Value = 1323
$ws.Range("L5").Value = 6817.5
$ws.Range("M5").Value = -1211
$ws.Range("N5").Value = -7041.5
$ws.Range("H8").Value = 90.5
$ws.Range("I8").Value = 90.5
$ws.Range("K8").Value = 271.5
$ws.Range("M8").Value = -132.5
$ws.Range("H122").Value = 2347.4695
$ws.Range("I122").Value = 700.13635
$ws.Range("J122").Value = 3689.7407
$ws.Range("K122").Value = 6301.22715
$ws.Range("L122").Value = 33207.6663
$ws.Range("M122").Value = -3851.22715
$ws.Range("N122").Value = -38107.6663
$ws.Range("H135").Value = 1404.9474
$ws.Range("I135").Value = 441
$ws.Range("J135").Value = 2272.5
$ws.Range("K135").Value = 3969
$ws.Range("L135").Value = 20452.5
$ws.Range("M135").Value = -1434
$ws.Range("N135").Value = -25522.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 16149.667
$ws.Range("J39").Value = 16149.667
$ws.Range("L39").Value = 16149.667
$ws.Range("N39").Value = -17213.667
$ws.Range("H80").Value = 2447.6667
$ws.Range("I80").Value = 2226.5
$ws.Range("K80").Value = 2226.5
$ws.Range("M80").Value = -1228.5
$ws.Range("H83").Value = 2447.6667
$ws.Range("I83").Value = 2226.5
$ws.Range("K83").Value = 11132.5
$ws.Range("M83").Value = -6140.5
$ws.Range("H132").Value = 2249.5356
$ws.Range("I132").Value = 1157.3125
$ws.Range("K132").Value = 3471.9375
$ws.Range("M132").Value = -941.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14185
$ws.Range("I132").Value = 17567.133
$ws.Range("K132").Value = 52701.399
$ws.Range("M132").Value = -50171.399

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 4450
$ws.Range("I19").Value = 2900
$ws.Range("J19").Value = 6000
$ws.Range("K19").Value = 2900
$ws.Range("L19").Value = 6000
$ws.Range("M19").Value = -2726
$ws.Range("N19").Value = -6348
$ws.Range("H81").Value = 1605.4166
$ws.Range("I81").Value = 1718.1111
$ws.Range("J81").Value = 1267.3334
$ws.Range("K81").Value = 3436.2222
$ws.Range("L81").Value = 2534.6668
$ws.Range("M81").Value = -2375.2222
$ws.Range("N81").Value = -4656.6668
$ws.Range("H84").Value = 1605.4166
$ws.Range("I84").Value = 1718.1111
$ws.Range("J84").Value = 1267.3334
$ws.Range("K84").Value = 17181.111
$ws.Range("L84").Value = 12673.334
$ws.Range("M84").Value = -11877.111
$ws.Range("N84").Value = -23281.334
$ws.Range("H122").Value = 4328.8887
$ws.Range("I122").Value = 2576.6667
$ws.Range("J122").Value = 7833.3335
$ws.Range("K122").Value = 7730.000100000001
$ws.Range("L122").Value = 23500.0005
$ws.Range("M122").Value = -5280.000100000001
$ws.Range("N122").Value = -28400.0005
$ws.Range("H132").Value = 1847.5853
$ws.Range("J132").Value = 2894.2856
$ws.Range("L132").Value = 8682.856800000001
$ws.Range("N132").Value = -13742.8568
